# Replace single-colon "air:" compartment labels with double-colon "air::" labels
# in column B (the only column containing these values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = $ws.Range("B1:B231")
$colB.Replace("air:", "air::") | Out-Null

# Update the active cell selection to match the saved worksheet view.
$ws.Range("F39").Select() | Out-Null
